$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.055200099945068
$ws.Range("B1").Value = 3.735349178314209
$ws.Range("C1").Value = 2.017575263977051
$ws.Range("D1").Value = 1.607144713401794
$ws.Range("E1").Value = 1.478622555732727
